$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.093609571456909
$ws.Range("B1").Value = 1.603262305259705
$ws.Range("C1").Value = 3.525516748428345
$ws.Range("D1").Value = 3.666451454162598
$ws.Range("E1").Value = 0.9699335098266602
